# Apply odds updates to Sheet1 as described by the commit's diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("AI2").Value = 23

# Row 3
$ws.Range("AY3").Value = 41

# Row 5
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 13
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 3.75
$ws.Range("Q5").Value = 1.85
$ws.Range("R5").Value = 2

# Row 6
$ws.Range("P6").Value = 4

# Row 8
$ws.Range("G8").Value = 1.9
$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 2.5
$ws.Range("K8").Value = 2.3
$ws.Range("N8").Value = 12
$ws.Range("S8").Value = 1.33
$ws.Range("T8").Value = 3.25
$ws.Range("U8").Value = 1.62
$ws.Range("V8").Value = 2.2
$ws.Range("AE8").Value = 12
$ws.Range("AT8").Value = 3.25
$ws.Range("AW8").Value = 6

# Row 9
$ws.Range("N9").Value = 13
$ws.Range("Q9").Value = 1.73
$ws.Range("R9").Value = 2.08

# Row 10
$ws.Range("H10").Value = 3.3
$ws.Range("K10").Value = 2.05
$ws.Range("O10").Value = 1.4
$ws.Range("P10").Value = 2.75
$ws.Range("Q10").Value = 2.25
$ws.Range("R10").Value = 1.62
$ws.Range("S10").Value = 1.5
$ws.Range("T10").Value = 2.5
$ws.Range("U10").Value = 2.1
$ws.Range("V10").Value = 1.67
$ws.Range("AC10").Value = 7.5
$ws.Range("AE10").Value = 19
$ws.Range("AH10").Value = 10
$ws.Range("AI10").Value = 21
$ws.Range("AM10").Value = 51
$ws.Range("AN10").Value = 3.6
$ws.Range("AO10").Value = 10
$ws.Range("AR10").Value = 67
$ws.Range("AT10").Value = 2.5
$ws.Range("AY10").Value = 41
$ws.Range("BA10").Value = 151
$ws.Range("BB10").Value = 351

# Row 11
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("Q11").Value = 2.3
$ws.Range("R11").Value = 1.6
